$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 8).Value = 0.1171128981563193
$ws.Cells.Item(3, 2).Value = 0.001751616726639724
$ws.Cells.Item(3, 3).Value = 0.0005886542447783803
$ws.Cells.Item(3, 4).Value = 4.72041620568
$ws.Cells.Item(3, 5).Value = 0.07148561940509333
$ws.Cells.Item(3, 6).Value = 0.0005978714841086201
$ws.Cells.Item(3, 7).Value = 0.002905361969170828
$ws.Cells.Item(3, 8).Value = 0.118864514882959
$ws.Cells.Item(4, 2).Value = 0.005976259980203255
$ws.Cells.Item(4, 3).Value = 0.0008752086631142932
$ws.Cells.Item(4, 4).Value = 9.070296014769472
$ws.Cells.Item(4, 5).Value = 0.1160658002708508
$ws.Cells.Item(4, 6).Value = 0.004260875680703753
$ws.Cells.Item(4, 7).Value = 0.007691644279702758
$ws.Cells.Item(4, 8).Value = 0.1230891581365225
$ws.Cells.Item(5, 2).Value = 0.01450988904263016
$ws.Cells.Item(5, 3).Value = 0.001869879344199833
$ws.Cells.Item(5, 4).Value = 12.08588329835214
$ws.Cells.Item(5, 5).Value = 0.116797703649839
$ws.Cells.Item(5, 6).Value = 0.01084497900738589
$ws.Cells.Item(5, 7).Value = 0.01817479907787442
$ws.Cells.Item(5, 8).Value = 0.1316227871989495
$ws.Cells.Item(6, 2).Value = 0.02361377326808725
$ws.Cells.Item(6, 3).Value = 0.00479548871852964
$ws.Cells.Item(6, 4).Value = 12.10583851656908
$ws.Cells.Item(6, 5).Value = 0.8828958137500963
$ws.Cells.Item(6, 6).Value = 0.01583627778076365
$ws.Cells.Item(6, 7).Value = 0.0289753680946028
$ws.Cells.Item(6, 8).Value = 0.1407266714244066
$ws.Cells.Item(7, 2).Value = 0.024533271278712
$ws.Cells.Item(7, 3).Value = 0.005990063738319108
$ws.Cells.Item(7, 4).Value = 10.64061017633164
$ws.Cells.Item(7, 5).Value = 0.6999744282950865
$ws.Cells.Item(7, 6).Value = 0.01428901431217685
$ws.Cells.Item(7, 7).Value = 0.03274618708407508
$ws.Cells.Item(7, 8).Value = 0.1416461694350313
$ws.Cells.Item(8, 2).Value = 0.02731791392427544
$ws.Cells.Item(8, 3).Value = 0.004659349870548872
$ws.Cells.Item(8, 4).Value = 11.13449745625414
$ws.Cells.Item(8, 5).Value = -0.08109131265847877
$ws.Cells.Item(8, 6).Value = 0.0187742668980547
$ws.Cells.Item(8, 7).Value = 0.03707169085692393
$ws.Cells.Item(8, 8).Value = 0.1444308120805947
$ws.Cells.Item(9, 2).Value = 0.02895797071699005
$ws.Cells.Item(9, 3).Value = 0.003142553185564316
$ws.Cells.Item(9, 4).Value = 9.764336544709577
$ws.Cells.Item(9, 5).Value = -0.1432627292233572
$ws.Cells.Item(9, 6).Value = 0.02601714830501916
$ws.Cells.Item(9, 7).Value = 0.03585942964113151
$ws.Cells.Item(9, 8).Value = 0.1460708688733094
$ws.Cells.Item(10, 2).Value = -0.1171128981563193
$ws.Cells.Item(10, 3).Value = 0.0004661576400478731
$ws.Cells.Item(10, 4).Value = -260.007564558441
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 6).Value = -0.1180265535176245
$ws.Cells.Item(10, 7).Value = -0.116199242795014
$ws.Cells.Item(11, 2).Value = -0.05914300480267276
$ws.Cells.Item(11, 3).Value = 0.0005233340776124656
$ws.Cells.Item(11, 4).Value = -113.6928165133035
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(11, 6).Value = -0.06016872429518282
$ws.Cells.Item(11, 7).Value = -0.05811728531016268
$ws.Cells.Item(11, 8).Value = 0.05796989335364654
$ws.Cells.Item(12, 2).Value = -0.05055369132759287
$ws.Cells.Item(12, 3).Value = 0.000509490071951381
$ws.Cells.Item(12, 4).Value = -101.0688135133905
$ws.Cells.Item(12, 5).Value = 0
$ws.Cells.Item(12, 6).Value = -0.05155227696881815
$ws.Cells.Item(12, 7).Value = -0.04955510568636758
$ws.Cells.Item(12, 8).Value = 0.06655920682872643
$ws.Cells.Item(13, 2).Value = -0.04238585770763247
$ws.Cells.Item(13, 3).Value = 0.0005055313894263784
$ws.Cells.Item(13, 4).Value = -84.00008364807073
$ws.Cells.Item(13, 5).Value = [double]"3.864917428580268e-304"
$ws.Cells.Item(13, 6).Value = -0.04337668443873938
$ws.Cells.Item(13, 7).Value = -0.04139503097652557
$ws.Cells.Item(13, 8).Value = 0.07472704044868683
$ws.Cells.Item(14, 2).Value = -0.03715318570259091
$ws.Cells.Item(14, 3).Value = 0.0005013035115643763
$ws.Cells.Item(14, 4).Value = -72.56134769261372
$ws.Cells.Item(14, 5).Value = [double]"7.264147222171206e-198"
$ws.Cells.Item(14, 6).Value = -0.03813572593228513
$ws.Cells.Item(14, 7).Value = -0.03617064547289668
$ws.Cells.Item(14, 8).Value = 0.07995971245372839
$ws.Cells.Item(15, 2).Value = -0.03166955018955077
$ws.Cells.Item(15, 3).Value = 0.0004928681619996101
$ws.Cells.Item(15, 4).Value = -62.618241039164
$ws.Cells.Item(15, 5).Value = [double]"1.016287705386315e-94"
$ws.Cells.Item(15, 6).Value = -0.03263555740145004
$ws.Cells.Item(15, 7).Value = -0.03070354297765149
$ws.Cells.Item(15, 8).Value = 0.08544334796676853
$ws.Cells.Item(16, 2).Value = -0.02944736298860311
$ws.Cells.Item(16, 3).Value = 0.0004746102597525732
$ws.Cells.Item(16, 4).Value = -58.83123800757404
$ws.Cells.Item(16, 5).Value = [double]"6.655003808821432e-28"
$ws.Cells.Item(16, 6).Value = -0.03037758527920681
$ws.Cells.Item(16, 7).Value = -0.02851714069799943
$ws.Cells.Item(16, 8).Value = 0.08766553516771619
$ws.Cells.Item(17, 2).Value = -0.02655517660578525
$ws.Cells.Item(17, 3).Value = 0.0004843724248204951
$ws.Cells.Item(17, 4).Value = -53.50460004612993
$ws.Cells.Item(17, 5).Value = [double]"1.872241760688606e-59"
$ws.Cells.Item(17, 6).Value = -0.02750453243099215
$ws.Cells.Item(17, 7).Value = -0.02560582078057835
$ws.Cells.Item(17, 8).Value = 0.09055772155053404
$ws.Cells.Item(18, 2).Value = -0.02412900377462856
$ws.Cells.Item(18, 3).Value = 0.0004883729551201607
$ws.Cells.Item(18, 4).Value = -48.46407167942537
$ws.Cells.Item(18, 5).Value = [double]"1.095091326729609e-14"
$ws.Cells.Item(18, 6).Value = -0.02508620048047725
$ws.Cells.Item(18, 7).Value = -0.02317180706877986
$ws.Cells.Item(18, 8).Value = 0.09298389438169075
$ws.Cells.Item(19, 2).Value = -0.0188972891002038
$ws.Cells.Item(19, 3).Value = 0.0004789421010117074
$ws.Cells.Item(19, 4).Value = -38.83378889115956
$ws.Cells.Item(19, 5).Value = 0.004069961672764923
$ws.Cells.Item(19, 6).Value = -0.01983600163436389
$ws.Cells.Item(19, 7).Value = -0.01795857656604369
$ws.Cells.Item(19, 8).Value = 0.0982156090561155
$ws.Cells.Item(20, 2).Value = -0.01593499849413125
$ws.Cells.Item(20, 3).Value = 0.000490662270992007
$ws.Cells.Item(20, 4).Value = -32.75281312191834
$ws.Cells.Item(20, 5).Value = 0.001804839781867012
$ws.Cells.Item(20, 6).Value = -0.0168966822541038
$ws.Cells.Item(20, 7).Value = -0.01497331473415869
$ws.Cells.Item(20, 8).Value = 0.1011778996621881
$ws.Cells.Item(21, 2).Value = -0.0116172847166318
$ws.Cells.Item(21, 3).Value = 0.0004928123625478579
$ws.Cells.Item(21, 4).Value = -24.89990881517644
$ws.Cells.Item(21, 5).Value = [double]"1.090203038055371e-09"
$ws.Cells.Item(21, 6).Value = -0.01258318254692239
$ws.Cells.Item(21, 7).Value = -0.0106513868863412
$ws.Cells.Item(21, 8).Value = 0.1054956134396875
$ws.Cells.Item(22, 2).Value = -0.009045888562279947
$ws.Cells.Item(22, 3).Value = 0.0004840828551545576
$ws.Cells.Item(22, 4).Value = -20.01838321270301
$ws.Cells.Item(22, 5).Value = [double]"1.434485411480741e-05"
$ws.Cells.Item(22, 6).Value = -0.009994676818258151
$ws.Cells.Item(22, 7).Value = -0.00809710030630174
$ws.Cells.Item(22, 8).Value = 0.1080670095940394
$ws.Cells.Item(23, 2).Value = -0.006839539416897422
$ws.Cells.Item(23, 3).Value = 0.0004781462641888565
$ws.Cells.Item(23, 4).Value = -15.66232366301672
$ws.Cells.Item(23, 5).Value = 0.03075118075715774
$ws.Cells.Item(23, 6).Value = -0.007776692137867619
$ws.Cells.Item(23, 7).Value = -0.005902386695927224
$ws.Cells.Item(23, 8).Value = 0.1102733587394219
$ws.Cells.Item(24, 2).Value = -0.003845330171996173
$ws.Cells.Item(24, 3).Value = 0.0004724253566438683
$ws.Cells.Item(24, 4).Value = -9.186877028688576
$ws.Cells.Item(24, 5).Value = 0.09880880719529377
$ws.Cells.Item(24, 6).Value = -0.00477127012937318
$ws.Cells.Item(24, 7).Value = -0.002919390214619164
$ws.Cells.Item(24, 8).Value = 0.1132675679843231
$ws.Cells.Item(25, 2).Value = -0.0017349681679533
$ws.Cells.Item(25, 3).Value = 0.0004531672323743885
$ws.Cells.Item(25, 4).Value = -4.734285572795315
$ws.Cells.Item(25, 5).Value = 0.07161111809629701
$ws.Cells.Item(25, 6).Value = -0.002623162764777247
$ws.Cells.Item(25, 7).Value = -0.0008467735711293535
$ws.Cells.Item(25, 8).Value = 0.115377929988366
$ws.Cells.Item(26, 2).Value = 0.06609685353543178
$ws.Cells.Item(26, 3).Value = 0.002642900245841702
$ws.Cells.Item(26, 4).Value = 33.95510765128341
$ws.Cells.Item(26, 5).Value = 0.0854206001918808
$ws.Cells.Item(26, 6).Value = 0.06091684773179267
$ws.Cells.Item(26, 7).Value = 0.07127685933907085
$ws.Cells.Item(26, 8).Value = 0.1832097516917511
